$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The "_GoBack" bookmark currently sits between "docum" and "ent"
#    in the "Replace the .txt ..." paragraph. Word moves this bookmark
#    to track the location of the most recent edit, so remove it from
#    its old spot - it will be re-created around the paragraph that
#    is actually edited below.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2) Merge ". Name the docum" + "ent anything you want to. " into a
#    single run with the combined text (removing the bookmark split
#    the two runs apart; replacing the whole span collapses them back
#    into one run again).
# ------------------------------------------------------------------
$mergeRange = $d.Content
$mergeRange.Find.Execute(". Name the docum", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $mergeRange.Start
$afterRange = $d.Content
$afterRange.Find.Execute("ent anything you want to. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeEnd = $afterRange.End

$wholeRange = $d.Range($mergeStart, $mergeEnd)
$wholeRange.Delete()
$wholeRange.InsertAfter(". Name the document anything you want to. ")

# ------------------------------------------------------------------
# 3) Remove bold formatting from the "Access to command prompt and
#    file explorer" paragraph (both the paragraph mark and the run).
# ------------------------------------------------------------------
$target = $d.Paragraphs(4)
$target.Range.Bold = 0

# ------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark spanning the paragraph that was
#    just edited, including its paragraph mark.
# ------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $target.Range) | Out-Null
